# Update B2:B451 (total energy use, scenario 3) with corrected model output values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 450,1
$arr[0,0] = 2138804981.158982
$arr[1,0] = 2139526041.86859
$arr[2,0] = 2142097914.686991
$arr[3,0] = 2143942338.260793
$arr[4,0] = 2145905736.913566
$arr[5,0] = 2147978882.529763
$arr[6,0] = 2150154224.376956
$arr[7,0] = 2152424702.733198
$arr[8,0] = 2154784410.039358
$arr[9,0] = 2157227945.907773
$arr[10,0] = 2159751222.042228
$arr[11,0] = 2162349955.807718
$arr[12,0] = 2165021243.927349
$arr[13,0] = 2167761864.231871
$arr[14,0] = 2170569453.038857
$arr[15,0] = 2173441791.101225
$arr[16,0] = 2176377184.191926
$arr[17,0] = 2179374122.756116
$arr[18,0] = 2182431511.866107
$arr[19,0] = 2185548650.309197
$arr[20,0] = 2188724510.197807
$arr[21,0] = 2191958323.450936
$arr[22,0] = 2195248978.190609
$arr[23,0] = 2198596799.092504
$arr[24,0] = 2202001589.231447
$arr[25,0] = 2205463490.502785
$arr[26,0] = 2208982371.227073
$arr[27,0] = 2212558413.586568
$arr[28,0] = 2216191652.774745
$arr[29,0] = 2219882561.703291
$arr[30,0] = 2223630946.461398
$arr[31,0] = 2227437800.992536
$arr[32,0] = 2231301418.624073
$arr[33,0] = 2235224973.213418
$arr[34,0] = 2239206789.579248
$arr[35,0] = 2243250111.474703
$arr[36,0] = 2247352396.291162
$arr[37,0] = 2251516295.356815
$arr[38,0] = 2255740800.341888
$arr[39,0] = 2260028805.54682
$arr[40,0] = 2264379568.94417
$arr[41,0] = 2268795069.831059
$arr[42,0] = 2273274337.160438
$arr[43,0] = 2277819259.891344
$arr[44,0] = 2282429786.712303
$arr[45,0] = 2287108143.925932
$arr[46,0] = 2291855055.130577
$arr[47,0] = 2296670125.132672
$arr[48,0] = 2301555053.005863
$arr[49,0] = 2306509970.40515
$arr[50,0] = 2311537170.028351
$arr[51,0] = 2316637704.817971
$arr[52,0] = 2321811153.48953
$arr[53,0] = 2327058246.167996
$arr[54,0] = 2332381328.200587
$arr[55,0] = 2337782073.85403
$arr[56,0] = 2343259558.214925
$arr[57,0] = 2348815399.136883
$arr[58,0] = 2354452422.471277
$arr[59,0] = 2360169555.79668
$arr[60,0] = 2365968395.442823
$arr[61,0] = 2371850821.89049
$arr[62,0] = 2377818612.072618
$arr[63,0] = 2383870640.39083
$arr[64,0] = 2390009770.285069
$arr[65,0] = 2396237834.388535
$arr[66,0] = 2402553321.830643
$arr[67,0] = 2408959926.720694
$arr[68,0] = 2415459273.310175
$arr[69,0] = 2422049740.218425
$arr[70,0] = 2428736437.768002
$arr[71,0] = 2435516804.447958
$arr[72,0] = 2442395301.818408
$arr[73,0] = 2449372227.895644
$arr[74,0] = 2456447835.425693
$arr[75,0] = 2463624972.627228
$arr[76,0] = 2470903880.050624
$arr[77,0] = 2478287761.977693
$arr[78,0] = 2485775968.665645
$arr[79,0] = 2493372032.247406
$arr[80,0] = 2501075718.939737
$arr[81,0] = 2508889332.332716
$arr[82,0] = 2516813310.345017
$arr[83,0] = 2524851276.436495
$arr[84,0] = 2533002824.300894
$arr[85,0] = 2541271753.286829
$arr[86,0] = 2549656438.740538
$arr[87,0] = 2558161895.84574
$arr[88,0] = 2566787222.717062
$arr[89,0] = 2575536232.659854
$arr[90,0] = 2584408001.475127
$arr[91,0] = 2593406691.744878
$arr[92,0] = 2602533076.907247
$arr[93,0] = 2611788346.711464
$arr[94,0] = 2621176617.28282
$arr[95,0] = 2630695896.646787
$arr[96,0] = 2640350250.011707
$arr[97,0] = 2650142364.751352
$arr[98,0] = 2660072896.860748
$arr[99,0] = 2670143798.40632
$arr[100,0] = 2680357948.407447
$arr[101,0] = 2690715105.088033
$arr[102,0] = 2701219306.561803
$arr[103,0] = 2711873049.420057
$arr[104,0] = 2722675200.549982
$arr[105,0] = 2733631570.315687
$arr[106,0] = 2744742255.984163
$arr[107,0] = 2756009735.887812
$arr[108,0] = 2767434804.146401
$arr[109,0] = 2779021878.581501
$arr[110,0] = 2790773134.910877
$arr[111,0] = 2802688418.192833
$arr[112,0] = 2814771564.640289
$arr[113,0] = 2827026587.839605
$arr[114,0] = 2839452795.94188
$arr[115,0] = 2852052187.014485
$arr[116,0] = 2864829675.25618
$arr[117,0] = 2877787695.977056
$arr[118,0] = 2890927333.839855
$arr[119,0] = 2904250068.774662
$arr[120,0] = 2917760158.826776
$arr[121,0] = 2931459553.53048
$arr[122,0] = 2945352710.094906
$arr[123,0] = 2959437228.402859
$arr[124,0] = 2973722137.247392
$arr[125,0] = 2988205531.816839
$arr[126,0] = 3002892936.090302
$arr[127,0] = 3017783451.976906
$arr[128,0] = 3032884253.115394
$arr[129,0] = 3048194950.387314
$arr[130,0] = 3063720610.448123
$arr[131,0] = 3079460870.256594
$arr[132,0] = 3095423294.419376
$arr[133,0] = 3111607028.031505
$arr[134,0] = 3128017320.979705
$arr[135,0] = 3144654804.190224
$arr[136,0] = 3161527177.342176
$arr[137,0] = 3178631307.987397
$arr[138,0] = 3195975582.971251
$arr[139,0] = 3213560766.0699
$arr[140,0] = 3231391704.877385
$arr[141,0] = 3249466740.547135
$arr[142,0] = 3267797404.102962
$arr[143,0] = 3286380693.974748
$arr[144,0] = 3305222521.362089
$arr[145,0] = 3324325355.453869
$arr[146,0] = 3343696030.509551
$arr[147,0] = 3363331806.199879
$arr[148,0] = 3383243253.998585
$arr[149,0] = 3403428636.946932
$arr[150,0] = 3423894646.986347
$arr[151,0] = 3444643836.415153
$arr[152,0] = 3465681267.692554
$arr[153,0] = 3487009280.933312
$arr[154,0] = 3508634544.213119
$arr[155,0] = 3530556081.013491
$arr[156,0] = 3552781184.174569
$arr[157,0] = 3575312609.986283
$arr[158,0] = 3598158127.590591
$arr[159,0] = 3621318691.211617
$arr[160,0] = 3644799838.185056
$arr[161,0] = 3668604807.04778
$arr[162,0] = 3692739696.391671
$arr[163,0] = 3717208618.857799
$arr[164,0] = 3742014301.826403
$arr[165,0] = 3767162409.444612
$arr[166,0] = 3792659895.78921
$arr[167,0] = 3818505905.528387
$arr[168,0] = 3844712216.438004
$arr[169,0] = 3871277772.242587
$arr[170,0] = 3898212640.375989
$arr[171,0] = 3925515148.94544
$arr[172,0] = 3953199312.960942
$arr[173,0] = 3981260249.368802
$arr[174,0] = 4009712807.551022
$arr[175,0] = 4038553334.945908
$arr[176,0] = 4067797290.805941
$arr[177,0] = 4097438541.270284
$arr[178,0] = 4127493779.688363
$arr[179,0] = 4157958936.954741
$arr[180,0] = 4188849133.287833
$arr[181,0] = 4220160432.087076
$arr[182,0] = 4251909402.761743
$arr[183,0] = 4284090274.032903
$arr[184,0] = 4316720740.710213
$arr[185,0] = 4349797672.842723
$arr[186,0] = 4383334101.613087
$arr[187,0] = 4417329931.840367
$arr[188,0] = 4451799247.885956
$arr[189,0] = 4486740543.991501
$arr[190,0] = 4522168254.925788
$arr[191,0] = 4558081632.611065
$arr[192,0] = 4594494324.483061
$arr[193,0] = 4631407312.677647
$arr[194,0] = 4668834015.067895
$arr[195,0] = 4706776844.459642
$arr[196,0] = 4745244131.795829
$arr[197,0] = 4784244551.964924
$arr[198,0] = 4823786103.464449
$arr[199,0] = 4863873086.459952
$arr[200,0] = 4904517928.598398
$arr[201,0] = 4945726500.568723
$arr[202,0] = 4987505446.350499
$arr[203,0] = 5029865199.270443
$arr[204,0] = 5072813361.398611
$arr[205,0] = 5116358169.515008
$arr[206,0] = 5160509471.734128
$arr[207,0] = 5205275245.721012
$arr[208,0] = 5250665139.983875
$arr[209,0] = 5296685945.423244
$arr[210,0] = 5343349923.180239
$arr[211,0] = 5390665758.456144
$arr[212,0] = 5438642170.930059
$arr[213,0] = 5487290321.551441
$arr[214,0] = 5536619348.078555
$arr[215,0] = 5586638711.593951
$arr[216,0] = 5637361407.398458
$arr[217,0] = 5688796118.59201
$arr[218,0] = 5740953178.780901
$arr[219,0] = 5793844903.003432
$arr[220,0] = 5847483169.822895
$arr[221,0] = 5901876105.423342
$arr[222,0] = 5957037982.571296
$arr[223,0] = 6012980543.61484
$arr[224,0] = 6069714345.286394
$arr[225,0] = 6127254326.255496
$arr[226,0] = 6185611014.122827
$arr[227,0] = 6244798943.813603
$arr[228,0] = 6304825529.740165
$arr[229,0] = 6365713936.896604
$arr[230,0] = 6427470357.718646
$arr[231,0] = 6490110075.312921
$arr[232,0] = 6553647784.661823
$arr[233,0] = 6618101004.753188
$arr[234,0] = 6683477867.191691
$arr[235,0] = 6749798687.168451
$arr[236,0] = 6817078651.030168
$arr[237,0] = 6885329150.255516
$arr[238,0] = 6954572233.517376
$arr[239,0] = 7024820873.574192
$arr[240,0] = 7096090891.704216
$arr[241,0] = 7168399492.147058
$arr[242,0] = 7241770272.071197
$arr[243,0] = 7316211050.427534
$arr[244,0] = 7391748737.630458
$arr[245,0] = 7468395787.193706
$arr[246,0] = 7546176555.243697
$arr[247,0] = 7625105305.304536
$arr[248,0] = 7705205717.111272
$arr[249,0] = 7786496590.43299
$arr[250,0] = 7868996238.212211
$arr[251,0] = 7952732583.170587
$arr[252,0] = 8037719354.285689
$arr[253,0] = 8123984348.915772
$arr[254,0] = 8211547587.165306
$arr[255,0] = 8300434981.789351
$arr[256,0] = 8390663631.652919
$arr[257,0] = 8482266349.54941
$arr[258,0] = 8575261374.695512
$arr[259,0] = 8669677535.694405
$arr[260,0] = 8765538237.943975
$arr[261,0] = 8862870283.251232
$arr[262,0] = 8961703531.850048
$arr[263,0] = 9062059884.901796
$arr[264,0] = 9163975246.232189
$arr[265,0] = 9267470037.710068
$arr[266,0] = 9372580025.519207
$arr[267,0] = 9479332641.519941
$arr[268,0] = 9587761608.529745
$arr[269,0] = 9697893045.30991
$arr[270,0] = 9809765562.721533
$arr[271,0] = 9923405869.558031
$arr[272,0] = 10038851062.59795
$arr[273,0] = 10156134218.44089
$arr[274,0] = 10275289757.55108
$arr[275,0] = 10396354095.74591
$arr[276,0] = 10519366732.61062
$arr[277,0] = 10644361757.26944
$arr[278,0] = 10771375718.69011
$arr[279,0] = 10900454494.12792
$arr[280,0] = 11031627485.80172
$arr[281,0] = 11164944490.64969
$arr[282,0] = 11300440420.99252
$arr[283,0] = 11438159511.56683
$arr[284,0] = 11578147978.45545
$arr[285,0] = 11720445165.23005
$arr[286,0] = 11865097638.55717
$arr[287,0] = 12012149093.07889
$arr[288,0] = 12161651023.31199
$arr[289,0] = 12313644624.23418
$arr[290,0] = 12468185936.69199
$arr[291,0] = 12625312914.57956
$arr[292,0] = 12785085232.99928
$arr[293,0] = 12947547664.21044
$arr[294,0] = 13112755073.65928
$arr[295,0] = 13280759844.04805
$arr[296,0] = 13451613324.18445
$arr[297,0] = 13625371335.16163
$arr[298,0] = 13802090652.03814
$arr[299,0] = 13981822338.97993
$arr[300,0] = 14164630886.03263
$arr[301,0] = 14316865095.61204
$arr[302,0] = 14471759563.98367
$arr[303,0] = 14629353058.59614
$arr[304,0] = 14789698148.11239
$arr[305,0] = 14952850799.16682
$arr[306,0] = 15118860598.54612
$arr[307,0] = 15287782488.72652
$arr[308,0] = 15459663962.49917
$arr[309,0] = 15634564012.07458
$arr[310,0] = 15812534433.84594
$arr[311,0] = 15993631216.24924
$arr[312,0] = 16177909952.05336
$arr[313,0] = 16365422703.38929
$arr[314,0] = 16556225362.0773
$arr[315,0] = 16750381664.72449
$arr[316,0] = 16947938180.14136
$arr[317,0] = 17148959036.83849
$arr[318,0] = 17353497735.02635
$arr[319,0] = 17561610093.2845
$arr[320,0] = 17773355746.92005
$arr[321,0] = 17988792053.13741
$arr[322,0] = 18207976198.66681
$arr[323,0] = 18430965310.22158
$arr[324,0] = 18657814102.19159
$arr[325,0] = 18888581611.40501
$arr[326,0] = 19123323975.74881
$arr[327,0] = 19362097059.14102
$arr[328,0] = 19604957387.10571
$arr[329,0] = 19851959653.66173
$arr[330,0] = 20103158949.27737
$arr[331,0] = 20358607673.41747
$arr[332,0] = 20618358957.44984
$arr[333,0] = 20882467960.69143
$arr[334,0] = 21150983209.36884
$arr[335,0] = 21423959160.83301
$arr[336,0] = 21701439516.65786
$arr[337,0] = 21983474311.46701
$arr[338,0] = 22270111473.59522
$arr[339,0] = 22561396075.4139
$arr[340,0] = 22857368510.54039
$arr[341,0] = 23158075438.95289
$arr[342,0] = 23463557755.88091
$arr[343,0] = 23773848413.77282
$arr[344,0] = 24088987277.4963
$arr[345,0] = 24409008197.39275
$arr[346,0] = 24733944459.18851
$arr[347,0] = 25063822817.59396
$arr[348,0] = 25398676771.10758
$arr[349,0] = 25738524268.39918
$arr[350,0] = 26083393169.86525
$arr[351,0] = 26433297929.07847
$arr[352,0] = 26788262932.11671
$arr[353,0] = 27148295814.74891
$arr[354,0] = 27513408511.0187
$arr[355,0] = 27883609297.33124
$arr[356,0] = 28252237558.24673
$arr[357,0] = 28625875518.05784
$arr[358,0] = 29004511798.47777
$arr[359,0] = 29388148814.02362
$arr[360,0] = 29776774308.21891
$arr[361,0] = 30170375878.98737
$arr[362,0] = 30568937836.48104
$arr[363,0] = 30972436463.10552
$arr[364,0] = 31380852880.38902
$arr[365,0] = 31794155446.35386
$arr[366,0] = 32212315008.31975
$arr[367,0] = 32635293746.5169
$arr[368,0] = 33063054568.24771
$arr[369,0] = 33495550472.87371
$arr[370,0] = 33932739554.90983
$arr[371,0] = 34311275615.94281
$arr[372,0] = 34693710047.01417
$arr[373,0] = 35079984701.03564
$arr[374,0] = 35470041607.92789
$arr[375,0] = 35863819421.95392
$arr[376,0] = 36261253198.65592
$arr[377,0] = 36661753917.9136
$arr[378,0] = 37059539412.53351
$arr[379,0] = 37439963178.847
$arr[380,0] = 37786029627.05816
$arr[381,0] = 38248501347.11165
$arr[382,0] = 38663081437.5509
$arr[383,0] = 39034438671.38519
$arr[384,0] = 39371987029.84505
$arr[385,0] = 39687037668.88828
$arr[386,0] = 39990788747.08349
$arr[387,0] = 40293129463.90365
$arr[388,0] = 40602127424.00042
$arr[389,0] = 40923928959.86397
$arr[390,0] = 41262907721.21414
$arr[391,0] = 41350907286.61491
$arr[392,0] = 41459517835.48065
$arr[393,0] = 41589311658.61166
$arr[394,0] = 41740150283.40096
$arr[395,0] = 41911389420.18561
$arr[396,0] = 42102059353.30245
$arr[397,0] = 42310968237.61712
$arr[398,0] = 42536817898.42505
$arr[399,0] = 42778256839.89666
$arr[400,0] = 43033944505.7524
$arr[401,0] = 43026516527.02568
$arr[402,0] = 43030275038.11356
$arr[403,0] = 43043815128.81816
$arr[404,0] = 43065846124.30353
$arr[405,0] = 43095196316.36259
$arr[406,0] = 43130801430.94424
$arr[407,0] = 43170978926.05836
$arr[408,0] = 43206927683.84193
$arr[409,0] = 43217959749.48016
$arr[410,0] = 43180804749.84083
$arr[411,0] = 43162704112.26219
$arr[412,0] = 43079106625.66346
$arr[413,0] = 42937736369.7973
$arr[414,0] = 42752189750.35725
$arr[415,0] = 42538008756.77411
$arr[416,0] = 42309999891.53535
$arr[417,0] = 42080788621.43114
$arr[418,0] = 41860271889.46875
$arr[419,0] = 41655635702.12622
$arr[420,0] = 41471652310.71738
$arr[421,0] = 41306590502.33797
$arr[422,0] = 41166223340.92298
$arr[423,0] = 41050638287.54202
$arr[424,0] = 40959099502.40549
$arr[425,0] = 41399554847.49578
$arr[426,0] = 41279468062.55145
$arr[427,0] = 41184150869.16489
$arr[428,0] = 41016668674.60599
$arr[429,0] = 40868264260.14008
$arr[430,0] = 40735366736.05183
$arr[431,0] = 40612883059.66239
$arr[432,0] = 40500785616.72161
$arr[433,0] = 40398335904.01894
$arr[434,0] = 40304158443.55113
$arr[435,0] = 40216144500.08947
$arr[436,0] = 40133370348.48991
$arr[437,0] = 40055503946.00752
$arr[438,0] = 39982047810.5573
$arr[439,0] = 39909836720.56779
$arr[440,0] = 39837415523.13437
$arr[441,0] = 39761056987.76433
$arr[442,0] = 39684117870.63231
$arr[443,0] = 39605223560.19102
$arr[444,0] = 39507960628.12733
$arr[445,0] = 39382520339.18952
$arr[446,0] = 39262673731.15862
$arr[447,0] = 39138942253.07367
$arr[448,0] = 39010867865.59238
$arr[449,0] = 38870841495.58385
$ws.Range("B2:B451").Value = $arr
